$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to convert an Excel serial date number (1899-12-30 epoch) to a .NET DateTime
function SerialToDate($serial) {
    $epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
    return $epoch.AddDays($serial)
}

$rows = @(
    @{ A = "LunaSensor";    B = "Sensor2"; C = 43381.313101851854; D = "Ekstra personale tilstede" },
    @{ A = "CarendoSensor"; B = "Sensor1"; C = 43380.862523148149; D = "Andet: Urolig borger" },
    @{ A = "CarendoSensor"; B = "Sensor3"; C = 43381.322974537034; D = "Ekstra personale tilstede" },
    @{ A = "CarendoSensor"; B = "Sensor3"; C = 43381.322974537034; D = "Ekstra personale tilstede" },
    @{ A = "LunaSensor";    B = "Sensor2"; C = 43381.313101851854; D = "Ekstra personale tilstede" }
)

$lastRow = $ws.UsedRange.Rows.Count
$startRow = $lastRow + 1
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B

    $cellC = $ws.Cells.Item($r, 3)
    # Copy formatting (date style) from an existing date cell, then set the value
    $ws.Cells.Item(2, 3).Copy($cellC)
    $cellC.Value = SerialToDate $row.C

    $ws.Cells.Item($r, 4).Value = $row.D
}
